# Add a new "Player Info" worksheet before the existing "ODI Batting" sheet,
# populate it with player metadata, and update the ODI Batting sheet's
# MATCH_CARD_LINK column to a simpler MATCH_CODE column.

$wb = $excel.ActiveWorkbook

# Existing sheet (will become the 2nd sheet).
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Remember the header style used on the existing sheet so the new sheet
# matches it.
$headerStyle = $battingSheet.Range("A1").Style

# Insert the new "Player Info" sheet before "ODI Batting".
$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

# NOTE: after Worksheets.Add(), the original $battingSheet handle now
# refers to the newly inserted sheet, so re-fetch the real "ODI Batting"
# sheet by name before touching it again.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Header row
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"
$infoSheet.Range("A1:D1").Style = $headerStyle

# Data row
$infoSheet.Range("A2").Value = "4427"
$infoSheet.Range("B2").Value = "Shane O Dowrich"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Does Not Bowl | Unknown"

# Update ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "4286"
